$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new staff member row (Staff ID, Name, Role, Gender, Age)
$ws.Range("A6").Value = "D999"
$ws.Range("B6").Value = "papa"
$ws.Range("C6").Value = "Doctor"
$ws.Range("D6").Value = "Male"
$ws.Range("E6").Value = 12

# Update the selection to reflect the newly added row (mirrors the author's
# UI state after entering the row: whole-row selection anchored at A6)
[void]$ws.Range("A6:XFD6").Select()
